$d = $word.ActiveDocument

# The document's sole paragraph (Title style) currently contains nothing
# but an empty "_GoBack" bookmark. Replace that bookmark with a run
# holding a single space character, inheriting the paragraph's run
# formatting (complex-script size 16pt / szCs 32).
$para = $d.Paragraphs.Item(1)
$rng = $para.Range
$rng.Delete()

$rng2 = $d.Paragraphs.Item(1).Range
$rng2.InsertAfter(" ")
$rng2.Font.SizeBi = 16
